$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.224.65"
$ws.Range("E2").Value = "  -2.10%  "

$ws.Range("D3").Value = "1.873.97"
$ws.Range("E3").Value = "  -1.48%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.002"
$ws.Range("E4").Value = "  -0.19%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "307.78"
$ws.Range("E5").Value = "  -1.62%  "

$ws.Range("E6").Value = "  -0.09%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5184"
$ws.Range("E7").Value = "  +3.35%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3760"
$ws.Range("E8").Value = "  -1.27%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07165"
$ws.Range("E9").Value = "  -1.52%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "20.86"
$ws.Range("E10").Value = "  +0.24%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.8863"
$ws.Range("E11").Value = "  -2.47%  "

$ws.Range("D12").Value = "1.888.11"
$ws.Range("E12").Value = "  -0.86%  "

$ws.Range("E13").Value = "  -0.69%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.347"
$ws.Range("E14").Value = "  -2.46%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "89.54"
$ws.Range("E15").Value = "  -1.89%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.002"
$ws.Range("E16").Value = "  -0.17%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.000008569"
$ws.Range("E17").Value = "  -1.53%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "14.23"
$ws.Range("E18").Value = "  -2.00%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "1.001"
$ws.Range("E19").Value = "  -0.15%  "

$ws.Range("D20").Value = "27.274.25"
$ws.Range("E20").Value = "  -2.04%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.049"

$ws.Range("D22").Value = "2.120.58"
$ws.Range("E22").Value = "  -0.99%  "

$ws.Range("E23").Value = "  -1.44%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.485"
$ws.Range("E24").Value = "  -1.59%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "151.80"

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.851"
$ws.Range("E26").Value = "  -0.50%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.195"
$ws.Range("E27").Value = "  -1.83%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "18.08"
$ws.Range("E28").Value = "  -1.53%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "113.10"
$ws.Range("E29").Value = "  -1.89%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.766"
$ws.Range("E30").Value = "  -2.91%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.718"
$ws.Range("E31").Value = "  +1.70%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.09058"
$ws.Range("E32").Value = "  +1.09%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.05186"
$ws.Range("E33").Value = "  -1.26%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.096"
$ws.Range("E34").Value = "  -3.40%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.7604"
$ws.Range("E35").Value = "  -0.76%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.185"
$ws.Range("E36").Value = "  -3.89%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.02052"
$ws.Range("E37").Value = "  -0.19%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.550"
$ws.Range("E38").Value = "  -0.05%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.041"
$ws.Range("E39").Value = "  +0.94%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.083"
$ws.Range("E40").Value = "  -1.19%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.5461"
$ws.Range("E41").Value = "  -1.49%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "6.704"
$ws.Range("E42").Value = "  -3.85%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "115.54"
$ws.Range("E43").Value = "  +4.08%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "8.572"
$ws.Range("E44").Value = "  +0.92%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.1492"
$ws.Range("E45").Value = "  -1.68%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.4714"
$ws.Range("E46").Value = "  -1.57%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "10.19"
$ws.Range("E47").Value = "  -3.74%  "

$ws.Range("E48").Value = "  -0.08%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.579"
$ws.Range("E49").Value = "  -3.09%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "65.21"
$ws.Range("E50").Value = "  -3.03%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "36.61"
$ws.Range("E51").Value = "  -1.24%  "
